$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update decimals (column G) and show_percent (column H) values for specific rows
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = $false

$ws.Range("G15").Value = 0

$ws.Range("G18").Value = 0

$ws.Range("H20").Value = $true

$ws.Range("G21").Value = 0

# Update the active selection to match the final state (G21 single cell)
$ws.Range("G21").Select()
